# feat: parse timestr with location timezone
#
# Rebuilds the start/end-time "chain" on the Activity sheet so that the
# activity segments run continuously from 2020-01-01 to 2020-10-01 in
# 5 steps (previously every row shared the same start/end values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity")

# New date boundaries used to subdivide the original 2020-01-01 .. 2020-10-01
# range into five consecutive activity windows.
$dates = @(
    "2020-01-01  05:00:00",
    "2020-03-05  05:00:00",
    "2020-07-10  05:00:00",
    "2020-08-10  05:00:00",
    "2020-09-10  05:00:00",
    "2020-10-01  05:00:00"
)

# Rows 2..6 correspond to dates[i-2] (start, column J) and dates[i-1] (end, column K)
for ($r = 2; $r -le 6; $r++) {
    $i = $r - 2
    $ws.Cells.Item($r, 10).Value = $dates[$i]
    $ws.Cells.Item($r, 11).Value = $dates[$i + 1]
}

# Update the saved selection on the Activity sheet.
$ws.Activate()
$ws.Range("K8").Select()
